$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("https://pglex.it/en/professionals/alessandra-arrighi/", "Alessandra Arrighi", "Senior Associate", "Pedersoli", "Italy", "-----", "-----", "aarrighi@pglex.it", "'39200220303051"),
    @("https://duncancotterill.com/people/amanda-bradley/", "Amanda Bradley", "Partner", "Duncan Cotterill", "New Zealand", "-----", "-----", "amanda.bradley@duncancotterill.com", "'64202720836204174"),
    @("https://ga-p.com/en/person/pablo-fernandez-cortijo/", "Pablo Fernández Cortijo", "Partner", "Gomez Acebo And Pombo", "Spain", "-----", "-----", "pfernandez@ga-p.com", "'34915829100"),
    @("https://ga-p.com/en/person/rodrigo-vilhena/", "Rodrigo Vilhena", "Senior Associate", "Gomez Acebo And Pombo", "Portugal", "-----", "-----", "rvcunha@ga-p.com", "'351213408600"),
    @("https://www.rdj.ie/our-people/padraic-brennan", "padraic brennan *****", "Partner", "RDJ", "Ireland", "-----", "-----", "padraic.brennan@rdj.ie", "'35391895313"),
    @("https://stikeman.com/en-ca/people/a/kevin-ackhurst", "Kevin Ackhurst", "Partner", "Stikeman Elliott", "Canada", "-----", "-----", "sabraham@stikeman.com", "'"),
    @("https://www.ogier.com/people/natalie-bell/", "Natalie Bell", "Partner", "Ogier", "the Cayman Islands", "-----", "Legal", "natalie.bell@ogier.com", "'13458151866"),
    @("https://www.ogier.com/people/christina-blacklaws/", "Christina Blacklaws", "Director", "Ogier", "Jersey", "-----", "Business Services Team", "christina.blacklaws@ogier.com", "'441534514000"),
    @("https://www.ogier.com/people/nicholas-burkill/", "Nicholas Burkill", "Partner", "Ogier", "the British Virgin Islands", "-----", "Legal", "nicholas.burkill@ogier.com", "'12848527372"),
    @("mailto:mirva.arvola@dittmar.fi", "Mirva Arvola", "Senior Associate", "Dittmar And Indrenius", "Finland", "-----", "-----", "mirva.arvola@dittmar.fi", "'358204020136203055"),
    @("https://wildlaw.ca/people/rory-cattanach", "rory *****", "'", "Wildeboer Dellelce", "Canada", "-----", "-----", "rory@wildlaw.ca", "'4163614766"),
    @("https://www.williamfry.com/people/colette-brady/", "Colette Brady", "Partner", "William Fry", "Ireland", "-----", "-----", "colette.brady@williamfry.com", "'35314896516"),
    @("https://www.williamfry.com/people/aisling-ogorman/", "Aisling O’gorman", "Partner", "William Fry", "England", "-----", "-----", "aisling.ogorman@williamfry.com", "'35316395167"),
    @("https://www.jgsa.pt/en/team/Pedro-Leao-Trigo/20/#!detail", "Pedro Leão Trigo", "Of Counsel", "JGSA", "Portugal", "-----", "Commercial and Corporate Law", "pleaotrigo@jgsa.pt", "'351213812690"),
    @("https://www.cassidylevy.com/person/isabel-fressynet/", "Isabel Fressynet", "Senior Associate", "Cassidy Levy Kent", "Belgium", "-----", "-----", "ifressynet@cassidylevy.com", "'32028831220"),
    @("https://www.cassidylevy.com/person/hugh-cheetham/", "Hugh Cheetham", "Counsel", "Cassidy Levy Kent", "Canada", "-----", "-----", "hcheetham@cassidylevy.com", "'16132620356"),
    @("https://www.hfw.com/people/anas-al-tarawneh/", "Anas Al-tarawneh", "Partner", "HFW", "Kuwait", "-----", "-----", "anas.altarawneh@hfw.com", "'97144230556"),
    @("https://www.hfw.com/people/sammy-beedan/", "Sammy Beedan Registered Foreign Lawyer", "Partner", "HFW", "Hong Kong", "-----", "-----", "sammy.beedan@hfw.com", "'85239837707"),
    @("https://www.hfw.com/people/simon-billing/", "Simon Billing", "Partner", "HFW", "Australia", "-----", "-----", "simon.billing@hfw.com", "'610894224716"),
    @("https://www.werksmans.com/our-people/chris-moraitis/", "Chris Moraitis", "Head", "Werksmans", "South Africa", "-----", "-----", "cmoraitis@werksmans.com", "'27115358271"),
    @("https://www.hna.ca/en/expert/gcarbonneau/", "Guy Carbonneau", "Partner", "HNA", "Canada", "-----", "-----", "gcarbonneau@hna.ca", "'14388026949"),
    @("https://www.magnussonlaw.com/people/ligita-biliune/", "Ligita Biliūnė", "Senior Associate", "Magnusson Law", "Lithuania", "-----", "Corporate and M&A", "ligita.kuliesaite@magnussonlaw.com", "'37052487388"),
    @("https://www.magnussonlaw.com/people/morvarid-dorkhan-nilsson/", "Morvarid Dorkhan Nilsson", "Partner", "Magnusson Law", "Sweden", "-----", "Public procurement", "morvarid.dorkhan@magnussonlaw.com", "'4684637500"),
    @("https://www.minterellison.co.nz/people/holly-hill", "Holly Hill", "Partner", "Minter Ellison RuddWatts", "New Zealand", "-----", "Property and real estate", "holly.hill@minterellison.co.nz", "'6444985173")
)

$rowIndex = 2
foreach ($rowData in $data) {
    $colIndex = 1
    foreach ($val in $rowData) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}
